# Applies the Shinryu_Profits scheduled-runner recompute: refreshed profit-margin
# figures (columns H:N) for the rows whose source/vendor prices moved.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Cells.Item(12, 8).Value = 110  # H12
$ws.Cells.Item(12, 9).Value = 113.333336  # I12
$ws.Cells.Item(12, 10).Value = 100  # J12
$ws.Cells.Item(12, 11).Value = 113.333336  # K12
$ws.Cells.Item(12, 12).Value = 100  # L12
$ws.Cells.Item(12, 13).Value = 56.666664  # M12
$ws.Cells.Item(12, 14).Value = -440  # N12
# Row 52
$ws.Cells.Item(52, 8).Value = 2000  # H52
$ws.Cells.Item(52, 9).Value = 2000  # I52
$ws.Cells.Item(52, 11).Value = 6000  # K52
$ws.Cells.Item(52, 13).Value = -5840  # M52
# Row 64
$ws.Cells.Item(64, 8).Value = 3369.3774  # H64
$ws.Cells.Item(64, 9).Value = 3213.2666  # I64
$ws.Cells.Item(64, 10).Value = 3573  # J64
$ws.Cells.Item(64, 11).Value = 3213.2666  # K64
$ws.Cells.Item(64, 12).Value = 3573  # L64
$ws.Cells.Item(64, 13).Value = -2965.2666  # M64
$ws.Cells.Item(64, 14).Value = -4069  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 3369.3774  # H67
$ws.Cells.Item(67, 9).Value = 3213.2666  # I67
$ws.Cells.Item(67, 10).Value = 3573  # J67
$ws.Cells.Item(67, 11).Value = 3213.2666  # K67
$ws.Cells.Item(67, 12).Value = 3573  # L67
$ws.Cells.Item(67, 13).Value = -2355.2666  # M67
$ws.Cells.Item(67, 14).Value = -5289  # N67
# Row 113
$ws.Cells.Item(113, 8).Value = 1947.0588  # H113
$ws.Cells.Item(113, 9).Value = 1550  # I113
$ws.Cells.Item(113, 10).Value = 2000  # J113
$ws.Cells.Item(113, 11).Value = 1550  # K113
$ws.Cells.Item(113, 12).Value = 2000  # L113
$ws.Cells.Item(113, 13).Value = 1704  # M113
$ws.Cells.Item(113, 14).Value = -8508  # N113
# Row 116
$ws.Cells.Item(116, 8).Value = 1878376.9  # H116
$ws.Cells.Item(116, 9).Value = 8548564  # I116
$ws.Cells.Item(116, 10).Value = 2386.75  # J116
$ws.Cells.Item(116, 11).Value = 8548564  # K116
$ws.Cells.Item(116, 12).Value = 2386.75  # L116
$ws.Cells.Item(116, 13).Value = -8545122  # M116
$ws.Cells.Item(116, 14).Value = -9270.75  # N116
# Row 132
$ws.Cells.Item(132, 8).Value = 3744.9355  # H132
$ws.Cells.Item(132, 9).Value = 4685.9546  # I132
$ws.Cells.Item(132, 10).Value = 1444.6666  # J132
$ws.Cells.Item(132, 11).Value = 14057.8638  # K132
$ws.Cells.Item(132, 12).Value = 4333.9998  # L132
$ws.Cells.Item(132, 13).Value = -11527.8638  # M132
$ws.Cells.Item(132, 14).Value = -9393.9998  # N132
# Row 137
$ws.Cells.Item(137, 8).Value = 3770.634  # H137
$ws.Cells.Item(137, 9).Value = 1181.0625  # I137
$ws.Cells.Item(137, 10).Value = 5427.96  # J137
$ws.Cells.Item(137, 11).Value = 3543.1875  # K137
$ws.Cells.Item(137, 12).Value = 16283.88  # L137
$ws.Cells.Item(137, 13).Value = -993.1875  # M137
$ws.Cells.Item(137, 14).Value = -21383.88  # N137
# Row 138
$ws.Cells.Item(138, 8).Value = 1867.8727  # H138
$ws.Cells.Item(138, 9).Value = 979.9524  # I138
$ws.Cells.Item(138, 10).Value = 2416.2942  # J138
$ws.Cells.Item(138, 11).Value = 2939.8572  # K138
$ws.Cells.Item(138, 12).Value = 7248.882599999999  # L138
$ws.Cells.Item(138, 13).Value = 2200.1428  # M138
$ws.Cells.Item(138, 14).Value = -17528.8826  # N138

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 1750  # H45
$ws.Cells.Item(45, 9).Value = 918  # I45
$ws.Cells.Item(45, 10).Value = 3414  # J45
$ws.Cells.Item(45, 11).Value = 918  # K45
$ws.Cells.Item(45, 12).Value = 3414  # L45
$ws.Cells.Item(45, 13).Value = -541  # M45
$ws.Cells.Item(45, 14).Value = -4168  # N45
# Row 61
$ws.Cells.Item(61, 8).Value = 1677.2972  # H61
$ws.Cells.Item(61, 9).Value = 1228.2424  # I61
$ws.Cells.Item(61, 11).Value = 1228.2424  # K61
$ws.Cells.Item(61, 13).Value = -1016.2424  # M61
# Row 122
$ws.Cells.Item(122, 8).Value = 1537.1666  # H122
$ws.Cells.Item(122, 9).Value = 1537.1666  # I122
$ws.Cells.Item(122, 10).Value = 0  # J122
$ws.Cells.Item(122, 11).Value = 4611.4998  # K122
$ws.Cells.Item(122, 12).Value = 0  # L122
$ws.Cells.Item(122, 13).Value = -2161.4998  # M122
$ws.Cells.Item(122, 14).ClearContents()  # N122
# Row 136
$ws.Cells.Item(136, 8).Value = 1677.2972  # H136
$ws.Cells.Item(136, 9).Value = 1228.2424  # I136
$ws.Cells.Item(136, 11).Value = 3684.7272  # K136
$ws.Cells.Item(136, 13).Value = -1134.7272  # M136

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2557  # H16
$ws.Cells.Item(16, 9).Value = 981.4286  # I16
$ws.Cells.Item(16, 10).Value = 6233.3335  # J16
$ws.Cells.Item(16, 11).Value = 981.4286  # K16
$ws.Cells.Item(16, 12).Value = 6233.3335  # L16
$ws.Cells.Item(16, 13).Value = -694.4286  # M16
$ws.Cells.Item(16, 14).Value = -6807.3335  # N16
# Row 58
$ws.Cells.Item(58, 8).Value = 1567.5745  # H58
$ws.Cells.Item(58, 9).Value = 934.0769  # I58
$ws.Cells.Item(58, 10).Value = 2351.9048  # J58
$ws.Cells.Item(58, 11).Value = 934.0769  # K58
$ws.Cells.Item(58, 12).Value = 2351.9048  # L58
$ws.Cells.Item(58, 13).Value = -731.0769  # M58
$ws.Cells.Item(58, 14).Value = -2757.9048  # N58
# Row 99
$ws.Cells.Item(99, 8).Value = 2010.85  # H99
$ws.Cells.Item(99, 9).Value = 1900.1666  # I99
$ws.Cells.Item(99, 10).Value = 3007  # J99
$ws.Cells.Item(99, 11).Value = 1900.1666  # K99
$ws.Cells.Item(99, 12).Value = 3007  # L99
$ws.Cells.Item(99, 13).Value = -402.1666  # M99
$ws.Cells.Item(99, 14).Value = -6003  # N99
# Row 107
$ws.Cells.Item(107, 8).Value = 1005.7  # H107
$ws.Cells.Item(107, 9).Value = 771.8333  # I107
$ws.Cells.Item(107, 10).Value = 1356.5  # J107
$ws.Cells.Item(107, 11).Value = 771.8333  # K107
$ws.Cells.Item(107, 12).Value = 1356.5  # L107
$ws.Cells.Item(107, 13).Value = 1148.1667  # M107
$ws.Cells.Item(107, 14).Value = -5196.5  # N107
# Row 113
$ws.Cells.Item(113, 8).Value = 2557  # H113
$ws.Cells.Item(113, 9).Value = 981.4286  # I113
$ws.Cells.Item(113, 10).Value = 6233.3335  # J113
$ws.Cells.Item(113, 11).Value = 981.4286  # K113
$ws.Cells.Item(113, 12).Value = 6233.3335  # L113
$ws.Cells.Item(113, 13).Value = 1188.5714  # M113
$ws.Cells.Item(113, 14).Value = -10573.3335  # N113
# Row 126
$ws.Cells.Item(126, 8).Value = 2010.85  # H126
$ws.Cells.Item(126, 9).Value = 1900.1666  # I126
$ws.Cells.Item(126, 10).Value = 3007  # J126
$ws.Cells.Item(126, 11).Value = 5700.4998  # K126
$ws.Cells.Item(126, 12).Value = 9021  # L126
$ws.Cells.Item(126, 13).Value = -3230.4998  # M126
$ws.Cells.Item(126, 14).Value = -13961  # N126
# Row 134
$ws.Cells.Item(134, 8).Value = 4584.625  # H134
$ws.Cells.Item(134, 9).Value = 2358.8572  # I134
$ws.Cells.Item(134, 10).Value = 6315.778  # J134
$ws.Cells.Item(134, 11).Value = 7076.571599999999  # K134
$ws.Cells.Item(134, 12).Value = 18947.334  # L134
$ws.Cells.Item(134, 13).Value = -4541.571599999999  # M134
$ws.Cells.Item(134, 14).Value = -24017.334  # N134
# Row 136
$ws.Cells.Item(136, 8).Value = 1567.5745  # H136
$ws.Cells.Item(136, 9).Value = 934.0769  # I136
$ws.Cells.Item(136, 10).Value = 2351.9048  # J136
$ws.Cells.Item(136, 11).Value = 2802.2307  # K136
$ws.Cells.Item(136, 12).Value = 7055.714399999999  # L136
$ws.Cells.Item(136, 13).Value = -252.2307000000001  # M136
$ws.Cells.Item(136, 14).Value = -12155.7144  # N136

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Cells.Item(109, 8).Value = 3549.3447  # H109
$ws.Cells.Item(109, 9).Value = 370.1111  # I109
$ws.Cells.Item(109, 11).Value = 1110.3333  # K109
$ws.Cells.Item(109, 13).Value = -70.33330000000001  # M109

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 2941.2307  # H80
$ws.Cells.Item(80, 9).Value = 2608.25  # I80
$ws.Cells.Item(80, 10).Value = 3474  # J80
$ws.Cells.Item(80, 11).Value = 2608.25  # K80
$ws.Cells.Item(80, 12).Value = 3474  # L80
$ws.Cells.Item(80, 13).Value = -1610.25  # M80
$ws.Cells.Item(80, 14).Value = -5470  # N80
# Row 83
$ws.Cells.Item(83, 8).Value = 2941.2307  # H83
$ws.Cells.Item(83, 9).Value = 2608.25  # I83
$ws.Cells.Item(83, 10).Value = 3474  # J83
$ws.Cells.Item(83, 11).Value = 13041.25  # K83
$ws.Cells.Item(83, 12).Value = 17370  # L83
$ws.Cells.Item(83, 13).Value = -8049.25  # M83
$ws.Cells.Item(83, 14).Value = -27354  # N83
# Row 102
$ws.Cells.Item(102, 8).Value = 944.6829  # H102
$ws.Cells.Item(102, 9).Value = 823.06665  # I102
$ws.Cells.Item(102, 10).Value = 1276.3636  # J102
$ws.Cells.Item(102, 11).Value = 823.06665  # K102
$ws.Cells.Item(102, 12).Value = 1276.3636  # L102
$ws.Cells.Item(102, 13).Value = 798.93335  # M102
$ws.Cells.Item(102, 14).Value = -4520.3636  # N102
# Row 132
$ws.Cells.Item(132, 8).Value = 2667.0789  # H132
$ws.Cells.Item(132, 9).Value = 2242.6897  # I132
$ws.Cells.Item(132, 10).Value = 4034.5557  # J132
$ws.Cells.Item(132, 11).Value = 6728.0691  # K132
$ws.Cells.Item(132, 12).Value = 12103.6671  # L132
$ws.Cells.Item(132, 13).Value = -4198.0691  # M132
$ws.Cells.Item(132, 14).Value = -17163.6671  # N132

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Cells.Item(93, 8).Value = 8281.429  # H93
$ws.Cells.Item(93, 9).Value = 17706  # I93
$ws.Cells.Item(93, 10).Value = 1213  # J93
$ws.Cells.Item(93, 11).Value = 17706  # K93
$ws.Cells.Item(93, 12).Value = 1213  # L93
$ws.Cells.Item(93, 13).Value = -16458  # M93
$ws.Cells.Item(93, 14).Value = -3709  # N93
# Row 111
$ws.Cells.Item(111, 8).Value = 0  # H111
$ws.Cells.Item(111, 10).Value = 0  # J111
$ws.Cells.Item(111, 12).Value = 0  # L111
$ws.Cells.Item(111, 14).ClearContents()  # N111
# Row 122
$ws.Cells.Item(122, 8).Value = 3099.2903  # H122
$ws.Cells.Item(122, 9).Value = 2834.0435  # I122
$ws.Cells.Item(122, 10).Value = 3861.875  # J122
$ws.Cells.Item(122, 11).Value = 8502.130500000001  # K122
$ws.Cells.Item(122, 12).Value = 11585.625  # L122
$ws.Cells.Item(122, 13).Value = -6052.130500000001  # M122
$ws.Cells.Item(122, 14).Value = -16485.625  # N122

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Cells.Item(136, 8).Value = 1498.3429  # H136
$ws.Cells.Item(136, 9).Value = 702.36365  # I136
$ws.Cells.Item(136, 10).Value = 2845.3845  # J136
$ws.Cells.Item(136, 11).Value = 2107.09095  # K136
$ws.Cells.Item(136, 12).Value = 8536.1535  # L136
$ws.Cells.Item(136, 13).Value = 442.9090500000002  # M136
$ws.Cells.Item(136, 14).Value = -13636.1535  # N136

Write-Output "All edits applied"